# "Generate Report for Handback" - refresh the localization-status report after
# a handback round-trip completed: the overall status text flips from
# "Ready for handoff" to "Handed back: in sync with en-US", and the per-language
# sheets (zh-cn / de-de) get their "Latest Target File" / "Latest Handback File"
# / "Latest Handback DateTime" columns populated for both rows.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Source .md hyperlink target (same file referenced by column A / row 2 on both
# language sheets) - the "Latest Target File" links back to the same source doc.
$targetFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb7d87cc95b9a666e0bc05fd617330a2d71b7920/e2e/0a48022b-3f22-4b64-95fa-057cae1d5fe7.md"
$targetFileName = "0a48022b-3f22-4b64-95fa-057cae1d5fe7.md"

# ---------------------------------------------------------------------------
# Overview sheet: both rows' Status columns (zh-cn=E, de-de=F) move to the new
# "handed back" wording.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# zh-cn sheet: Status + handback columns for rows 2 and 3.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$zhHandbackFile = "0a48022b-3f22-4b64-95fa-057cae1d5fe7.1eae068af0547b479fed30e3d824c7431b952718.zh-cn.xlf"
$zhHandbackDate = "2016-08-31 21:18:36"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetFileUrl, "", "", $targetFileName) | Out-Null
$wsZh.Range("J2").Value = $zhHandbackFile
$wsZh.Range("K2").Value = $zhHandbackDate

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $targetFileUrl, "", "", $targetFileName) | Out-Null
$wsZh.Range("J3").Value = $zhHandbackFile
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------------
# de-de sheet: Status + handback columns for rows 2 and 3.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$deHandbackFile = "0a48022b-3f22-4b64-95fa-057cae1d5fe7.1eae068af0547b479fed30e3d824c7431b952718.de-de.xlf"
$deHandbackDate = "2016-08-31 21:18:44"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetFileUrl, "", "", $targetFileName) | Out-Null
$wsDe.Range("J2").Value = $deHandbackFile
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $targetFileUrl, "", "", $targetFileName) | Out-Null
$wsDe.Range("J3").Value = $deHandbackFile
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1
